$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and 1h-volume-change (E) values.
# D-column values are written via a Text number format + ClearFormats()
# round-trip so numeric-looking strings (e.g. '5.70', '1.00') are kept as
# literal text (matching the source data) instead of being coerced into
# doubles by Excel's usual 'looks like a number' inference -- while still
# ending up on the default (unstyled) cell format, same as the original file.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.895.56'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.423.01'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '552.37'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.02'
$ws.Range('D6').ClearFormats()
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.92%  '
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.70'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  -1.84%  '
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.24'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.855.80'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.859.54'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.95%  '
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.425.85'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.33'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.55%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.40'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '330.90'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.27'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.70%  '
$ws.Range('E24').Value = '  +0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.75'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.41%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('E27').Value = '  +3.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0778'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.78'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.34'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.14'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.68'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('E35').Value = '  +4.21%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.63'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.411'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '313.80'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +8.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.69'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '139.24'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('E44').Value = '  +1.09%  '
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.54'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.69%  '
$ws.Range('E47').Value = '  +2.07%  '
$ws.Range('E48').Value = '  +0.90%  '
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.06'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.38%  '
